$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.060.25'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').Value = '3.506.81'
$ws.Range('E3').Value = '  -0.93%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.56'
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.57'
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  -1.39%  '
$ws.Range('D8').Value = '3.501.64'
$ws.Range('E8').Value = '  -0.97%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('E10').Value = '  -2.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.26'
$ws.Range('E11').Value = '  +7.13%  '
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '46.05'
$ws.Range('E13').Value = '  -2.64%  '
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('D15').Value = '4.073.89'
$ws.Range('E15').Value = '  -1.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.37'
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '614.62'
$ws.Range('E17').Value = '  -2.18%  '
$ws.Range('D18').Value = '3.504.86'
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('D19').Value = '70.081.65'
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('E20').Value = '  +1.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.54'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.878'
$ws.Range('E22').Value = '  -0.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.16'
$ws.Range('E23').Value = '  -8.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '98.76'
$ws.Range('E24').Value = '  +2.17%  '
$ws.Range('E25').Value = '  -2.56%  '
$ws.Range('E26').Value = '  -3.25%  '
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.57'
$ws.Range('E28').Value = '  -1.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.75'
$ws.Range('E29').Value = '  +0.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.03'
$ws.Range('E30').Value = '  -1.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.99'
$ws.Range('E31').Value = '  -3.35%  '
$ws.Range('E32').Value = '  -5.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '632.52'
$ws.Range('E33').Value = '  +10.82%  '
$ws.Range('E34').Value = '  -4.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.81'
$ws.Range('E35').Value = '  -2.50%  '
$ws.Range('E36').Value = '  -2.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.75'
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('E38').Value = '  +3.70%  '
$ws.Range('E39').Value = '  -4.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '56.78'
$ws.Range('E40').Value = '  -1.23%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.146'
$ws.Range('E41').Value = '  +1.96%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').Value = '3.355.92'
$ws.Range('E43').Value = '  +0.54%  '
$ws.Range('D44').Value = '0.0₃0734'
$ws.Range('E44').Value = '  +3.07%  '
$ws.Range('E45').Value = '  -5.22%  '
$ws.Range('E46').Value = '  -4.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '31.91'
$ws.Range('E47').Value = '  -3.58%  '
$ws.Range('E48').Value = '  -3.98%  '
$ws.Range('E49').Value = '  +0.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '133.38'
$ws.Range('E50').Value = '  -0.25%  '
